$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.797777
$ws.Range("H2").Value = 122.393331
$ws.Range("I2").Value = 0.2689231481273683
$ws.Range("J2").Value = 0.2689231481273683
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.806204333333334
$ws.Range("N2").Value = 14.418613
$ws.Range("O2").Value = 0.7287437301541012
$ws.Range("P2").Value = 0.7287437301541012
$ws.Range("Q2").Value = 196.082452607767
$ws.Range("R2").Value = 1764.742073469903
$ws.Range("S2").Value = 0.1959760580911223
$ws.Range("T2").Value = 0.1959760580911223

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.797777
$ws.Range("H3").Value = 122.393331
$ws.Range("I3").Value = 0.2689231481273683
$ws.Range("J3").Value = 0.2689231481273683
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.788987
$ws.Range("N3").Value = 5.366961
$ws.Range("O3").Value = 0.2712562698458988
$ws.Range("P3").Value = 0.2712562698458988
$ws.Range("Q3").Value = 72.986692681899
$ws.Range("R3").Value = 656.8802341370911
$ws.Range("S3").Value = 0.07294709003624604
$ws.Range("T3").Value = 0.07294709003624605

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 46.219831
$ws.Range("H4").Value = 138.659493
$ws.Range("I4").Value = 0.3046632285488233
$ws.Range("J4").Value = 0.3046632285488233
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.806204333333334
$ws.Range("N4").Value = 14.418613
$ws.Range("O4").Value = 0.7287437301541012
$ws.Range("P4").Value = 0.7287437301541012
$ws.Range("Q4").Value = 222.1419520381343
$ws.Range("R4").Value = 1999.277568343209
$ws.Range("S4").Value = 0.222021417613461
$ws.Range("T4").Value = 0.222021417613461

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 46.219831
$ws.Range("H5").Value = 138.659493
$ws.Range("I5").Value = 0.3046632285488233
$ws.Range("J5").Value = 0.3046632285488233
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.788987
$ws.Range("N5").Value = 5.366961
$ws.Range("O5").Value = 0.2712562698458988
$ws.Range("P5").Value = 0.2712562698458988
$ws.Range("Q5").Value = 82.686676801197
$ws.Range("R5").Value = 744.180091210773
$ws.Range("S5").Value = 0.08264181093536237
$ws.Range("T5").Value = 0.08264181093536237

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 64.69033266666666
$ws.Range("H6").Value = 194.070998
$ws.Range("I6").Value = 0.4264136233238083
$ws.Range("J6").Value = 0.4264136233238083
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.806204333333334
$ws.Range("N6").Value = 14.418613
$ws.Range("O6").Value = 0.7287437301541012
$ws.Range("P6").Value = 0.7287437301541012
$ws.Range("Q6").Value = 310.9149571873082
$ws.Range("R6").Value = 2798.234614685774
$ws.Range("S6").Value = 0.3107462544495179
$ws.Range("T6").Value = 0.3107462544495179

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 64.69033266666666
$ws.Range("H7").Value = 194.070998
$ws.Range("I7").Value = 0.4264136233238083
$ws.Range("J7").Value = 0.4264136233238083
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.788987
$ws.Range("N7").Value = 5.366961
$ws.Range("O7").Value = 0.2712562698458988
$ws.Range("P7").Value = 0.2712562698458988
$ws.Range("Q7").Value = 115.730164166342
$ws.Range("R7").Value = 1041.571477497078
$ws.Range("S7").Value = 0.1156673688742904
$ws.Range("T7").Value = 0.1156673688742904
